$wb = $excel.ActiveWorkbook

# --- Sheet 1 (cover/metadata sheet) ---
$ws1 = $wb.Worksheets.Item(1)

# Remove the duplicated "Contact" row (old row 11); this shifts every
# subsequent row up by one and yields the final A1:B20 dimension.
$ws1.Rows.Item(11).Delete()

# Version bump
$ws1.Range("B3").Value = "6.0.0"

# Updated publication date
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws1.Range("B9").Value = "Alvearie Team"

# The remaining "Contact" row becomes the new "Jurisdiction" row
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# --- Sheet 2 (elements table) ---
$ws2 = $wb.Worksheets.Item(2)

# Fill in the Short / Definition text for the root Extension element
$ws2.Range("K2").Value = "Employee Exempt Indicator"
$ws2.Range("L2").Value = "Indicator of whether the employee status is exempt or non-exempt"
